$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 headers: columns D..K get re-arranged ---
$ws.Range("D1").Value = "Xtick Locations"
$ws.Range("E1").Value = "XTick Labels"
$ws.Range("F1").Value = "Kink Location"
$ws.Range("G1").Value = "Protocol Names"
$ws.Range("H1").Value = "Plotting Order"
$ws.Range("I1").Value = "Colors"
$ws.Range("J1").Value = "LineStyles"
$ws.Range("K1").ClearContents()

# --- Row 2 ---
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = 4
$ws.Range("G2").Value = "English Scaling"
$ws.Range("H2").Value = 1
$ws.Range("I2").Value = "blue"
$ws.Range("J2").Value = "-"
$ws.Range("K2").ClearContents()

# --- Row 3 ---
$ws.Range("E3").ClearContents()
$ws.Range("F3").ClearContents()
$ws.Range("G3").Value = "Long English Scaling"
$ws.Range("H3").Value = 3
$ws.Range("I3").Value = "red"
$ws.Range("J3").Value = "-"
$ws.Range("K3").ClearContents()

# --- Row 4 ---
$ws.Range("E4").ClearContents()
$ws.Range("F4").ClearContents()
$ws.Range("G4").Value = "Long Nonsense Scaling"
$ws.Range("H4").Value = 4
$ws.Range("I4").Value = "magenta"
$ws.Range("J4").Value = "-"
$ws.Range("K4").ClearContents()

# --- Row 5 ---
$ws.Range("E5").ClearContents()
$ws.Range("F5").ClearContents()
$ws.Range("G5").Value = "Nonsense Scaling"
$ws.Range("H5").Value = 2
$ws.Range("I5").Value = "cyan"
$ws.Range("J5").Value = "-"
$ws.Range("K5").ClearContents()

# --- Selection moves to F3 ---
$ws.Range("F3").Select()
